# Add a new "Talk: Data wrangling" section to the day-2 overview slide,
# right before the existing "Let's do stuff" bullet, mirroring the
# Talk/sub-bullets/blank-line pattern already used elsewhere in the deck.

$p = $ppt.ActivePresentation

# Locate the slide that contains the day 2 overview ("Talk: Clean code and
# unit testing" ... "Let's do stuff").
$s = $null
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $candidate = $p.Slides.Item($i)
    for ($j = 1; $j -le $candidate.Shapes.Count; $j++) {
        $shape = $candidate.Shapes.Item($j)
        if ($shape.HasTextFrame -and $shape.TextFrame.HasText) {
            if ($shape.TextFrame.TextRange.Text -like "*Clean code and unit testing*") {
                $s = $candidate
            }
        }
    }
}

$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

# Find the paragraph that currently reads "Let's do stuff" (top level, i.e.
# IndentLevel 1) - the new section must be inserted right before it.
# (TextRange.Paragraphs(i,1).Text includes the trailing paragraph-mark CR,
# so trim it before comparing.)
$targetIndex = -1
for ($i = 1; $i -le $tr.Paragraphs().Count; $i++) {
    $para = $tr.Paragraphs($i, 1)
    if ($para.IndentLevel -eq 1 -and $para.Text.TrimEnd("`r") -eq "Let's do stuff") {
        $targetIndex = $i
    }
}

$targetPara = $tr.Paragraphs($targetIndex, 1)

# Insert the four new paragraphs immediately before it:
#   Talk: Data wrangling      (level 1 / lvl 0)
#   What is it?                (level 2 / lvl 1)
#   Why do we need it?         (level 2 / lvl 1)
#   <empty>                    (level 2 / lvl 1)
$null = $targetPara.InsertBefore("Talk: Data wrangling`rWhat is it?`rWhy do we need it?`r`r")

# The three new sub-paragraphs need to be demoted to the second outline
# level to match the surrounding bullets.
$tr.Paragraphs($targetIndex + 1, 1).IndentLevel = 2
$tr.Paragraphs($targetIndex + 2, 1).IndentLevel = 2
$tr.Paragraphs($targetIndex + 3, 1).IndentLevel = 2
